# Intermediate code commit and updated the task sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New task rows (16-23), continuing the numbering/format of the existing table.
$tasks = @(
    @{ Row = 17; No = 16; Task = "Save Password" },
    @{ Row = 18; No = 17; Task = "Login --> Texboxes width increase" },
    @{ Row = 19; No = 18; Task = "Texboxes content place holder" },
    @{ Row = 20; No = 19; Task = "Form authentication. Redirect to login in case other url hit without login" },
    @{ Row = 21; No = 20; Task = "Toastr for message" },
    @{ Row = 22; No = 21; Task = "After login failed, no message shown to user." },
    @{ Row = 23; No = 22; Task = "Add User/ Manage Users" },
    @{ Row = 24; No = 23; Task = "Left Side Menu, Role based." }
)

foreach ($t in $tasks) {
    $r = $t.Row
    $ws.Cells.Item($r, 1).Value = $t.No
    $ws.Cells.Item($r, 2).Value = $t.Task
    $ws.Cells.Item($r, 5).Value = "Yet to decide"
    $ws.Cells.Item($r, 5).WrapText = $true
}

# Row 20 (long "Form authentication..." task) wraps onto two lines like the
# other multi-line entries in the sheet.
$ws.Rows.Item(20).RowHeight = 30

# Restore the on-screen view: scrolled down with row 10 at the top and D20
# selected (where the last edit happened).
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("D20").Select()
